$d = $word.ActiveDocument

# Helper: split the single character at [pos, pos+1) into its own run by
# re-typing it and toggling a character property (forces the engine to
# stop coalescing it with the neighbouring run of identical formatting).
function Split-CharRun($doc, $pos, $newChar) {
    $posEnd = $pos + 1
    $r = $doc.Range($pos, $posEnd)
    $r.Text = $newChar
    $r2 = $doc.Range($pos, $posEnd)
    $r2.Bold = 1
    $r2.Bold = 0
}

# The document contains three occurrences of the misspelled package name
# "mgcb" (should be "mgcv", the R package for GAMs). Only the first two
# (in the "Parte non parametrica" and "Semplifica e alloca al centroide"
# paragraphs) need to be fixed; the third ("Non complicato con mgcb...")
# is left untouched.

$searchStart = 0
for ($occurrence = 1; $occurrence -le 2; $occurrence++) {
    $docEnd = $d.Content.End
    $searchRange = $d.Range($searchStart, $docEnd)
    $found = $searchRange.Find.Execute("mgcb", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        break
    }

    $matchStart = $searchRange.Start
    $matchEnd = $searchRange.End

    # Replace the final "b" with "v" - this lands as its own run, matching
    # a user selecting just that character and retyping it.
    $lastCharPos = $matchEnd - 1
    Split-CharRun $d $lastCharPos "v"

    # In the first occurrence, the run that used to be
    # " di base sullo spazio (" also gets split into a lone space run and
    # "di base sullo spazio (" - mirror that same single-character retype.
    if ($occurrence -eq 1) {
        $afterWord = $matchEnd
        $afterWordEnd = $afterWord + 1
        $spaceRange = $d.Range($afterWord, $afterWordEnd)
        if ($spaceRange.Text -eq " ") {
            Split-CharRun $d $afterWord " "
        }
    }

    $searchStart = $matchEnd
}
